$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update session year
$ws.Range("G1").Value = "SESSION 2025"

# Update "Réalisations en milieu professionnel en cours de seconde année" entries
$ws.Range("A28").Value = "Gestion des services Active Directory (AD)"
$ws.Range("A29").Value = "Création et gestion des stratégies de groupe (GPO) "
$ws.Range("A30").Value = " Configuration et durcissement des switchs Cisco et HP"
$ws.Range("A32").Value = "Surveillance des performances des machines virtuelles sur ESXi"

# Update view: scroll position and active cell selection
$ws.Range("A1").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("G1").Select()
